$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel
# stores them as literal text (matching the source price-ticker feed)
# instead of coercing to a Double.
$textCells = @('D5','D6','D8','D10','D14','D16','D19','D21','D23','D24','D25','D29','D30','D31','D32','D39','D41','D43','D45','D46','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.218.36'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '1.644.27'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '217.16'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').Value = '0.515'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').Value = '20.00'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '1.873.91'
$ws.Range('D13').Value = '1.642.86'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('E15').Value = '  +3.24%  '
$ws.Range('D16').Value = '67.40'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '27.203.61'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('D19').Value = '218.92'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').Value = '6.85'
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('E22').Value = '  +5.42%  '
$ws.Range('D23').Value = '4.41'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('D24').Value = '9.19'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '147.70'
$ws.Range('E25').Value = '  +1.52%  '
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '15.75'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').Value = '3.37'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('E34').Value = '  +1.70%  '
$ws.Range('D35').Value = '1.261.26'
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').Value = '0.852'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D41').Value = '0.808'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('E42').Value = '  +6.38%  '
$ws.Range('D43').Value = '5.30'
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('D44').Value = '1.784.25'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '61.80'
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('D46').Value = '91.74'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.64'
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0973'
$ws.Range('E51').Value = '  +0.35%  '

Write-Host "Applied cryptos update"